$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values in B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values in B2:E2
$ws.Range("B2").Value = 511.961010133125
$ws.Range("C2").Value = 431.257796713125
$ws.Range("D2").Value = 511.961010133125
$ws.Range("E2").Value = 431.257796713125

# Update row 3 (STR) values in B3:E3
$ws.Range("B3").Value = 519.52693639124993
$ws.Range("C3").Value = 431.257796713125
$ws.Range("D3").Value = 519.52693639124993
$ws.Range("E3").Value = 431.257796713125

# Update the selected range to reflect the new active selection B1:E3
$ws.Range("B1:E3").Select()
